$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.220.23'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '1.990.66'
$ws.Range('E3').Value = '  +5.97%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '324.56'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5095'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4118'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +4.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08674'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +5.38%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.129'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.23%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '42.79'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.52%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '24.34'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +3.16%  '
$ws.Range('D13').Value = '1.998.62'
$ws.Range('E13').Value = '  +6.17%  '
$ws.Range('E14').Value = '  +2.84%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.380'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.50%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.9996'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '93.88'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.31%  '
$ws.Range('E18').Value = '  +2.11%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06550'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.30%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.75'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +3.77%  '
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.067'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.84%  '
$ws.Range('D23').Value = '30.291.92'
$ws.Range('E23').Value = '  +0.66%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.56'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.58%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.198'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.26%  '
$ws.Range('D26').Value = '2.226.78'
$ws.Range('E26').Value = '  +6.16%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.44'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +5.84%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '162.92'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.355'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +4.71%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '130.43'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.34%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.127'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +5.09%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1050'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.47%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.050'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.88%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.818'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +3.31%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.307'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +11.27%  '
$ws.Range('E36').Value = '  +2.20%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.375'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.50%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06515'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.47%  '
$ws.Range('E39').Value = '  +2.55%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.910'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +4.76%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6567'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +4.21%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.78'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +4.30%  '
$ws.Range('E43').Value = '  +0.70%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.54'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.64%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6097'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.21%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.191'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +4.66%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.657'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '124.21'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.64%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.222'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.07%  '
$ws.Range('E50').Value = '  +2.34%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06864'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.85%  '
